$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.575126262626263
$ws.Range("C2").Value = 0.0164141414141414
$ws.Range("D2").Value = 0.0599747474747475
$ws.Range("E2").Value = 0.945707070707071
$ws.Range("F2").Value = 0.0113636363636364
$ws.Range("G2").Value = 0.943181818181818
$ws.Range("H2").Value = 0.0138888888888889
$ws.Range("I2").Value = 0.72979797979798
$ws.Range("J2").Value = 0.0416666666666667
$ws.Range("K2").Value = 0.0391414141414141
$ws.Range("L2").Value = 0.0366161616161616
$ws.Range("M2").Value = 0.840277777777778
$ws.Range("N2").Value = 0.0101010101010101
$ws.Range("P2").Value = 0.0328282828282828
$ws.Range("Q2").Value = 0.928030303030303
$ws.Range("R2").Value = 0.00631313131313131
$ws.Range("S2").Value = 0.00378787878787879
$ws.Range("U2").Value = 0.00631313131313131
$ws.Range("W2").Value = 0.0113636363636364
$ws.Range("X2").Value = 0.00757575757575758
$ws.Range("B3").Value = 0.0656565656565657
$ws.Range("C3").Value = 0.347853535353535
$ws.Range("D3").Value = 0.841540404040404
$ws.Range("E3").Value = 0.0441919191919192
$ws.Range("F3").Value = 0.00126262626262626
$ws.Range("H3").Value = 0.0694444444444444
$ws.Range("I3").Value = 0.0366161616161616
$ws.Range("J3").Value = 0.130681818181818
$ws.Range("K3").Value = 0.920454545454545
$ws.Range("L3").Value = 0.95959595959596
$ws.Range("M3").Value = 0.109217171717172
$ws.Range("N3").Value = 0.919823232323232
$ws.Range("O3").Value = 0.0151515151515152
$ws.Range("P3").Value = 0.00126262626262626
$ws.Range("R3").Value = 0.852272727272727
$ws.Range("S3").Value = 0.983585858585859
$ws.Range("T3").Value = 0.0549242424242424
$ws.Range("U3").Value = 0.00757575757575758
$ws.Range("V3").Value = 0.0214646464646465
$ws.Range("W3").Value = 0.0328282828282828
$ws.Range("X3").Value = 0.00126262626262626
$ws.Range("B4").Value = 0.35290404040404
$ws.Range("C4").Value = 0.0252525252525253
$ws.Range("D4").Value = 0.00378787878787879
$ws.Range("E4").Value = 0.00378787878787879
$ws.Range("F4").Value = 0.945707070707071
$ws.Range("G4").Value = 0.053030303030303
$ws.Range("H4").Value = 0.00252525252525253
$ws.Range("I4").Value = 0.0195707070707071
$ws.Range("J4").Value = 0.0536616161616162
$ws.Range("K4").Value = 0.0378787878787879
$ws.Range("L4").Value = 0.00126262626262626
$ws.Range("M4").Value = 0.00631313131313131
$ws.Range("N4").Value = 0.00378787878787879
$ws.Range("P4").Value = 0.963383838383838
$ws.Range("Q4").Value = 0.00378787878787879
$ws.Range("R4").Value = 0.133838383838384
$ws.Range("U4").Value = 0.00252525252525253
$ws.Range("V4").Value = 0.053030303030303
$ws.Range("W4").Value = 0.953282828282828
$ws.Range("X4").Value = 0.954545454545455
$ws.Range("B5").Value = 0.00631313131313131
$ws.Range("C5").Value = 0.610479797979798
$ws.Range("D5").Value = 0.0921717171717172
$ws.Range("E5").Value = 0.00505050505050505
$ws.Range("F5").Value = 0.0416666666666667
$ws.Range("G5").Value = 0.00378787878787879
$ws.Range("H5").Value = 0.914141414141414
$ws.Range("I5").Value = 0.214015151515152
$ws.Range("J5").Value = 0.773989898989899
$ws.Range("K5").Value = 0.00126262626262626
$ws.Range("L5").Value = 0.00252525252525253
$ws.Range("M5").Value = 0.0441919191919192
$ws.Range("N5").Value = 0.0662878787878788
$ws.Range("O5").Value = 0.984848484848485
$ws.Range("P5").Value = 0.00252525252525253
$ws.Range("Q5").Value = 0.0681818181818182
$ws.Range("R5").Value = 0.00757575757575758
$ws.Range("S5").Value = 0.0126262626262626
$ws.Range("T5").Value = 0.945075757575758
$ws.Range("U5").Value = 0.983585858585859
$ws.Range("V5").Value = 0.92550505050505
$ws.Range("W5").Value = 0.00252525252525253
$ws.Range("X5").Value = 0.0366161616161616
